$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Linked List")

# Row 10 - LRU Cache
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "08/27/2025"
$ws.Range("A10").Style = "Normal"
$ws.Range("E10").Value = "Somewhat"
$ws.Range("G10").Value = "Yes"

# Row 11 - Merge K Sorted Lists
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "08/30/2025"
$ws.Range("A11").Style = "Normal"
$ws.Range("E11").Value = "No"
$ws.Range("G11").Value = "Yes"

# Row 12 - Reverse Nodes In K Group
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "08/30/2025"
$ws.Range("A12").Style = "Normal"
$ws.Range("E12").Value = "No"
$ws.Range("G12").Value = "Yes"
